$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("B2").Value = -0.2230573992474752
$ws.Range("C2").Value = 0.2854886207816306
$ws.Range("D2").Value = 0.127389058624141
$ws.Range("E2").Value = 0.3569160386199267
$ws.Range("F2").Value = 0.2891476032052339
$ws.Range("B3").Value = -0.2042933947234856
$ws.Range("C3").Value = 0.2861682746732524
$ws.Range("D3").Value = 0.1727834463265767
$ws.Range("E3").Value = 0.4156722823650582
$ws.Range("F3").Value = 0.3767868006342249
$ws.Range("B4").Value = -0.2126892945567086
$ws.Range("C4").Value = 0.3122101537896099
$ws.Range("D4").Value = 0.1603748835476225
$ws.Range("E4").Value = 0.4004683302679782
$ws.Range("F4").Value = 0.3544083123310924
$ws.Range("B5").Value = -0.1809481785306808
$ws.Range("C5").Value = 0.267150876214921
$ws.Range("D5").Value = 0.08966209893550295
$ws.Range("E5").Value = 0.2994363019667171
$ws.Range("F5").Value = 0.2502235823900796
$ws.Range("B6").Value = -0.1763213443630235
$ws.Range("C6").Value = 0.2237613264222619
$ws.Range("D6").Value = 0.09028189153338712
$ws.Range("E6").Value = 0.3004694519138129
$ws.Range("F6").Value = 0.2564559200924167
$ws.Range("B7").Value = -0.182531442650734
$ws.Range("C7").Value = 0.2539746169011773
$ws.Range("D7").Value = 0.08014903705961754
$ws.Range("E7").Value = 0.2831060526721701
$ws.Range("F7").Value = 0.2295326190139252
$ws.Range("B8").Value = -0.1031535539421949
$ws.Range("C8").Value = 0.1812455504833794
$ws.Range("D8").Value = 0.05114288788033396
$ws.Range("E8").Value = 0.2261479336194208
$ws.Range("F8").Value = 0.2204601520168994
$ws.Range("B9").Value = -0.135926032552076
$ws.Range("C9").Value = 0.1991159763380844
$ws.Range("D9").Value = 0.04509109148225762
$ws.Range("E9").Value = 0.2123466304942407
$ws.Range("F9").Value = 0.1998069261446269
$ws.Range("B10").Value = -0.3880619836346724
$ws.Range("C10").Value = 0.3880619836346724
$ws.Range("D10").Value = 0.1505921031424768
$ws.Range("E10").Value = 0.3880619836346724
